$p = $ppt.ActivePresentation

# --- 1. Slide 16: change the table's style id (GUID) -----------------------
# Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") has a 2-column table in
# shape 3; point it at the built-in table style instead of the custom one.
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{16C5D46E-62C3-4E07-8254-8E25FD516F44}")

# --- 2. Re-colour the presentation theme (Integral -> Office) --------------
# Swap the deck's colour scheme from the "Integral" theme colours to the
# default "Office" theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink), in that order.
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Item(1).RGB  = 0x000000   # dk1
$themeColors.Item(2).RGB  = 0xFFFFFF   # lt1
$themeColors.Item(3).RGB  = 0x6A5444   # dk2
$themeColors.Item(4).RGB  = 0xE6E6E7   # lt2
$themeColors.Item(5).RGB  = 0xD59B5B   # accent1
$themeColors.Item(6).RGB  = 0x317DED   # accent2
$themeColors.Item(7).RGB  = 0xA5A5A5   # accent3
$themeColors.Item(8).RGB  = 0x00C0FF   # accent4
$themeColors.Item(9).RGB  = 0xC47244   # accent5
$themeColors.Item(10).RGB = 0x47AD70   # accent6
$themeColors.Item(11).RGB = 0xC16305   # hlink
$themeColors.Item(12).RGB = 0x724F95   # folHlink
